$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, pushing existing rows 31-72 down to 32-73.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with a new data record (copy of the
# previous row 31 record, but with an updated date and volume).
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value = "La Araucanía"
$ws.Cells.Item(31, 4).Value = 44629
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100108
$ws.Cells.Item(31, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(31, 9).Value = 100108004
$ws.Cells.Item(31, 10).Value = "Papaya"
$ws.Cells.Item(31, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 50
$ws.Cells.Item(31, 14).Value = 20000
$ws.Cells.Item(31, 15).Value = 20000
$ws.Cells.Item(31, 16).Value = 20000
$ws.Cells.Item(31, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(31, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(31, 19).Value = 2000
$ws.Cells.Item(31, 20).Value = 10
